$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.661166191101074
$ws.Range("B1").Value = 6.496612071990967
$ws.Range("C1").Value = 5.432275772094727
$ws.Range("D1").Value = 6.476631164550781
$ws.Range("E1").Value = 3.97675609588623
